$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F2").Value = "The candidate will answer questions on creative principles and storytelling in AI-generated text prompts to assess understanding.`nK9: Creative principles and storytelling in AI generated text prompts (MED-MED-3004-1.1)"
$ws.Range("F3").Value = "The candidate will respond to questions on basic AI algorithms and models used in script/text generation to demonstrate knowledge.`nK10: Basic AI algorithms and models used in script/text generation (MED-MED-3004-1.1)"
$ws.Range("F4").Value = "The candidate will incorporate unique creative principles and storytelling elements into AI-generated story ideas during a practical exam.`nA8: Incorporate unique creative principles and storytelling elements into AI-generated story ideas to avoid generic replication (MED-MED-3004-1.1)"
$ws.Range("F5").Value = "The candidate will use AI-generated text techniques and methodologies to develop script elements in a practical exam.`nA5: Use AI-generated text techniques and methodologies to develop script elements (MED-MED-3004-1.1)"
$ws.Range("F6").Value = "The candidate will respond to a series of short answer questions related to AI-generated script ideation techniques for world-building, storyline and character development, scenario iterations.`nK1: AI-generated script ideation techniques for world-building, storyline and character development, scenario iterations (MED-MED-3004-1.1)"
$ws.Range("F7").Value = "The candidate will identify key terms and themes for input into Gen AI platform to ideate and iterate story elements in a practical exam.`nA7: Identify key terms and themes for input into Gen AI platform to ideate and iterate story elements for incorporation into text prompts (MED-MED-3004-1.1)"
$ws.Range("F8").Value = "The candidate will answer questions on Gen AI tool limitations and solutions to assess their knowledge.`nK3: Gen AI tool limitations and solutions (MED-MED-3004-1.1)"
$ws.Range("F9").Value = "The candidate will apply and adjust relevant prompts in the generative process to improve iterations during a practical exam.`nA6: Apply and adjust relevant prompts in the generative process to improve iterations (MED-MED-3004-1.1)"
$ws.Range("F10").Value = "The candidate will respond to questions on basic NLP techniques and tools relevant to AI text generation to demonstrate knowledge.`nK8: Basic Natural Language Processing (NLP) techniques and tools relevant to AI Text generation (MED-MED-3004-1.1)"
$ws.Range("F11").Value = "The candidate will utilize NLP techniques and tools to enhance the quality and effectiveness of AI-generated story elements in a practical exam.`nA2: Utilise NLP techniques and tools to enhance the quality and effectiveness of AI-generated story elements (MED-MED-3004-1.1)"
$ws.Range("F12").Value = "The candidate will respond to questions about AI tools for script analysis and market research to demonstrate knowledge.`nK6: AI tools for script analysis and market research (MED-MED-3004-1.1)"
$ws.Range("F13").Value = "The candidate will filter and script-edit Gen-AI output in a practical exam to demonstrate their skills.`nA4: Filter and Script-edit Gen-AI output (MED-MED-3004-1.1)"
$ws.Range("F14").Value = "The candidate will answer questions on Gen AI and ethics awareness to assess their understanding.`nK7: Gen AI and Ethics awareness (MED-MED-3004-1.1)"
$ws.Range("F15").Value = "The candidate will answer questions on best practices to minimize plagiarism risk to show understanding.`nK5: Best Practices to minimise plagiarism risk (MED-MED-3004-1.1)"
$ws.Range("F16").Value = "The candidate will apply ethical considerations in the selection of prompts and reference usage during a practical exam.`nA3: Apply ethical considerations in the selection of prompts and reference usage (MED-MED-3004-1.1)"
$ws.Range("F17").Value = "The candidate will explain awareness and correction of bias in LLMs training data reproduced in output through written responses.`nK4: Awareness and correction of bias in LLMs training data reproduced in output (MED-MED-3004-1.1)"
$ws.Range("F18").Value = "The candidate will perform a practical exercise on analyzing AI output for bias and taking corrective steps and submit their analysis.`nA1: Analyse AI output for bias and taking corrective steps (MED-MED-3004-1.1)"
$ws.Range("F19").Value = "The candidate will respond to questions about copyright law covering Gen AI usage to demonstrate understanding.`nK2: Copyright law covering Gen AI usage (MED-MED-3004-1.1)"
$ws.Range("F20").Value = "The candidate will identify copyright risk in Gen-AI production and avoid copyright infringement in a practical exam.`nA9: Identify copyright risk in Gen-AI production and avoid copyright infringement (MED-MED-3004-1.1)"
